$d = $word.ActiveDocument
$LB = [char]11

function Replace-WithBreaks($searchText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    $rng.Text = $newText
}

# 1. Objetivos paragraph
$old1 = "1. Entendimento da relação entre a termodinâmica de soluções e os diagramas de fases.2. Domínio da leitura de diagramas unários, binários e ternários (configuração do sistema para um estado termodinâmico, leitura de composições de fases e cálculo de suas quantidades).3. Relacionamento entre microestruturas e diagramas de fases.4. Entendimento da seqüência de eventos que ocorrem no processo de solidificação em equilíbrio e fora de equilíbrio."
$new1 = "1. Entendimento da relação entre a termodinâmica de soluções e os diagramas de fases." + $LB + `
        "2. Domínio da leitura de diagramas unários, binários e ternários (configuração do sistema para um estado termodinâmico, leitura de composições de fases e cálculo de suas quantidades)." + $LB + `
        "3. Relacionamento entre microestruturas e diagramas de fases." + $LB + `
        "4. Entendimento da seqüência de eventos que ocorrem no processo de solidificação em equilíbrio e fora de equilíbrio."
Replace-WithBreaks $old1 $new1

# 2. Programa resumido paragraph
$old2 = "A. Introdução; teoria básica de equilíbrio de fases;B. Sistemas unários; C. Sistemas binários; D. Sistemas ternários; E. Cálculo termodinâmico de diagramas de fases; F. Trabalho Prático."
$new2 = "A. Introdução; teoria básica de equilíbrio de fases;" + $LB + `
        "B. Sistemas unários; " + $LB + `
        "C. Sistemas binários; " + $LB + `
        "D. Sistemas ternários; " + $LB + `
        "E. Cálculo termodinâmico de diagramas de fases; " + $LB + `
        "F. Trabalho Prático."
Replace-WithBreaks $old2 $new2

# 3. Programa paragraph
$old3 = "1. Introdução; revisão da termodinâmica de soluções; teoria básica de equilíbrio de fases; curvas de energia livre versus composição; regra das fases; 2. Sistemas unários, equilíbrios bi-, mono- e invariantes;3. Sistemas binários isomorfos; a regra da alavanca; solidificação em equilíbrio e fora de equilíbrio; mínimos e máximos; 4. Sistemas eutéticos binários; solidificação e microetruturas de ligas hipoeutéticas, eutéticas e hipereutéticas; solidificação unidirecional com eutéticos; casos limites de eutéticos; 5. Sistemas eutetóides binários; solidificação e microetruturas de ligas hipoeutetóides, eutetói-des e hipereutetóides; o sistema Fe-C; 6. Sistemas monotéticos; sistemas monotetóides; sistemas metatéticos; transformações congruentes; 7. Sistemas peritéticos binários; resfriamento em equilíbrio e fora do equilíbrio de ligas peritéticas; sistemas peritetóides binários; sistemas sintéticos binários; 8. Sistemas ternários isomorfos; o triângulo de Gibbs; seções isotérmicas; projeções liquidus; seções verticais; máximos e mínimos; resfriamento em equilíbrio; 9. Equilíbrio ternário de três fases; regra da alavanca em campos trifásicos; resfriamento em equilíbrio; 10. Equilíbrio ternário de quatro fases: equilíbrio de classe I; equilíbrio de classe II e equilíbrio de classe III; 11. Transformações congruentes em sistemas ternários; sistemas ternários complexos; 12. Cálculo termodinâmico de diagramas de fases; Trabalho prático."
$new3 = "1. Introdução; revisão da termodinâmica de soluções; teoria básica de equilíbrio de fases; curvas de energia livre versus composição; regra das fases; " + $LB + `
        "2. Sistemas unários, equilíbrios bi-, mono- e invariantes;" + $LB + `
        "3. Sistemas binários isomorfos; a regra da alavanca; solidificação em equilíbrio e fora de equilíbrio; mínimos e máximos; " + $LB + `
        "4. Sistemas eutéticos binários; solidificação e microetruturas de ligas hipoeutéticas, eutéticas e hipereutéticas; solidificação unidirecional com eutéticos; casos limites de eutéticos; " + $LB + `
        "5. Sistemas eutetóides binários; solidificação e microetruturas de ligas hipoeutetóides, eutetói-des e hipereutetóides; o sistema Fe-C; " + $LB + `
        "6. Sistemas monotéticos; sistemas monotetóides; sistemas metatéticos; transformações congruentes; " + $LB + `
        "7. Sistemas peritéticos binários; resfriamento em equilíbrio e fora do equilíbrio de ligas peritéticas; sistemas peritetóides binários; sistemas sintéticos binários; " + $LB + `
        "8. Sistemas ternários isomorfos; o triângulo de Gibbs; seções isotérmicas; projeções liquidus; seções verticais; máximos e mínimos; resfriamento em equilíbrio; " + $LB + `
        "9. Equilíbrio ternário de três fases; regra da alavanca em campos trifásicos; resfriamento em equilíbrio; 10. Equilíbrio ternário de quatro fases: equilíbrio de classe I; equilíbrio de classe II e equilíbrio de classe III; " + $LB + `
        "11. Transformações congruentes em sistemas ternários; sistemas ternários complexos; " + $LB + `
        "12. Cálculo termodinâmico de diagramas de fases; Trabalho prático."
Replace-WithBreaks $old3 $new3

# 4. Critério (inside Avaliação paragraph) - only the one run's text, leave the trailing <w:br/> alone
$old4 = "As avaliações individuais, a participação nas resoluções dos exercícios e repostas aos questionários assim como a condução do trabalho prático e a apresentação dos resultados nas formas oral e escrita serão agrupadas em duas notas (N1 e N2) que comporão a nota final (NF). O critério para cálculo da nota final é: NF = (N1+ N2)/2Serão aprovados os alunos com NF ≥ 5,0Serão reprovados os alunos com NF < 3,0"
$new4 = "As avaliações individuais, a participação nas resoluções dos exercícios e repostas aos questionários assim como a condução do trabalho prático e a apresentação dos resultados nas formas oral e escrita serão agrupadas em duas notas (N1 e N2) que comporão a nota final (NF). O critério para cálculo da nota final é: " + $LB + `
        "NF = (N1+ N2)/2" + $LB + `
        "Serão aprovados os alunos com NF ≥ 5,0" + $LB + `
        "Serão reprovados os alunos com NF < 3,0"
Replace-WithBreaks $old4 $new4

# 5. Bibliografia paragraph
$old5 = "01. Gordon, P. Principles of Phase Diagrams in Materials Systems, McGraw-Hill, 1968.02. Rhines, F. N. Phase Diagrams in Metallurgy: Their Development and Applications, McGraw-Hill, 1956.03. Prince, A. Alloy Phase Equilibria, Elsevier, 1966.04. Massalski, T. B. Binary Alloys Phase Diagrams, ASM, Metals Park, Ohio, 1990.05. Alloy Phase Diagrams, ASM Handbook, Volume 3, ASM, Metals Park, Ohio, 1992.06. Hansen, M. Constitution of Binary Alloys, McGraw-Hill, 1958.07. Elliot, R. P. Constitution of Binary Alloys: First Supplement, McGraw-Hill, 1965.08. Shunk, F. A. Constitution of Binary Alloys: Second Supplement, McGraw-Hill, 1969.09. Levin, E. M. Phase Diagram for Ceramists, The American Ceramic Society, 1964.10. Rudman, P. S. Phase Stability in Metals and Alloys, McGraw-Hill, 1967.11. Kaufman, L. Computer Calculation of Phase Diagrams with Special Reference to Refractory Metals, Academic Press.12. Hack, K. The SGTE Casebook - Thermodynamics at Work. The Institut of Metals, London,6.13. Hillert, M. Phase Equilibria, Phase Diagrams and Phase Transformations. Cambridge University Press, Cambridge, 1998.14. Thermocalc version M manuals: User Guide and Examples, ThermoCalc AB, Stockholm, 1997."""
$new5 = "01. Gordon, P. Principles of Phase Diagrams in Materials Systems, McGraw-Hill, 1968." + $LB + `
        "02. Rhines, F. N. Phase Diagrams in Metallurgy: Their Development and Applications, McGraw-Hill, 1956." + $LB + `
        "03. Prince, A. Alloy Phase Equilibria, Elsevier, 1966." + $LB + `
        "04. Massalski, T. B. Binary Alloys Phase Diagrams, ASM, Metals Park, Ohio, 1990." + $LB + `
        "05. Alloy Phase Diagrams, ASM Handbook, Volume 3, ASM, Metals Park, Ohio, 1992." + $LB + `
        "06. Hansen, M. Constitution of Binary Alloys, McGraw-Hill, 1958." + $LB + `
        "07. Elliot, R. P. Constitution of Binary Alloys: First Supplement, McGraw-Hill, 1965." + $LB + `
        "08. Shunk, F. A. Constitution of Binary Alloys: Second Supplement, McGraw-Hill, 1969." + $LB + `
        "09. Levin, E. M. Phase Diagram for Ceramists, The American Ceramic Society, 1964." + $LB + `
        "10. Rudman, P. S. Phase Stability in Metals and Alloys, McGraw-Hill, 1967." + $LB + `
        "11. Kaufman, L. Computer Calculation of Phase Diagrams with Special Reference to Refractory Metals, Academic Press." + $LB + `
        "12. Hack, K. The SGTE Casebook - Thermodynamics at Work. The Institut of Metals, London,6." + $LB + `
        "13. Hillert, M. Phase Equilibria, Phase Diagrams and Phase Transformations. Cambridge University Press, Cambridge, 1998." + $LB + `
        "14. Thermocalc version M manuals: User Guide and Examples, ThermoCalc AB, Stockholm, 1997."""
Replace-WithBreaks $old5 $new5
